$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New id/variant values for each row (id recomputed directly from the
# variant text, lower-cased, with no preferred-variant grouping / no
# levenshtein-distance matching applied anymore).
$ws.Range("B2").Value = "#lammert"
$ws.Range("C2").Value = "Lammert"

$ws.Range("B3").Value = "#karek"
$ws.Range("C3").Value = "Karek"

$ws.Range("B4").Value = "#karel"
$ws.Range("C4").Value = "Karel"

$ws.Range("B5").Value = "#jeronimo"
$ws.Range("C5").Value = "Jeronimo"

$ws.Range("B6").Value = "#franzyn"
$ws.Range("C6").Value = "Franzyn"

$ws.Range("B7").Value = "#isabel"
$ws.Range("C7").Value = "Isabel"

$ws.Range("B8").Value = "#jan"
$ws.Range("C8").Value = "Jan"

$ws.Range("B9").Value = "#lubeert"
$ws.Range("C9").Value = "Lubeert"

$ws.Range("B10").Value = "#koenhert"
$ws.Range("C10").Value = "Koenhert"

$ws.Range("B11").Value = "#lambert"
$ws.Range("C11").Value = "Lambert"

$ws.Range("B12").Value = "#fransyn"
$ws.Range("C12").Value = "Fransyn"

$ws.Range("B13").Value = "#lubbert"
$ws.Range("C13").Value = "Lubbert"

$ws.Range("B14").Value = "#jsabel"
$ws.Range("C14").Value = "Jsabel"

$ws.Range("B15").Value = "#izabel"
$ws.Range("C15").Value = "Izabel"

$ws.Range("B16").Value = "#hoogadel"
$ws.Range("C16").Value = "Hoogadel"

# is_prefered column is no longer populated on export.
$ws.Range("D2").Value = ""
$ws.Range("D3").Value = ""
$ws.Range("D4").Value = ""
$ws.Range("D5").Value = ""
$ws.Range("D6").Value = ""
$ws.Range("D7").Value = ""
$ws.Range("D8").Value = ""
$ws.Range("D9").Value = ""
$ws.Range("D10").Value = ""
